$d = $word.ActiveDocument

# Suppression du champ "rating" dans la liste des champs de la table
# "movies" : il est remplacé par la table intermédiaire dédiée
# ("movies_users_ratings"), donc ce champ n'a plus lieu d'être sur movies.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "rating (en int(11)*") {
        $p.Range.Delete()
        break
    }
}

Write-Output "done"
